$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the default (unstyled) look of data cells so forcing text-type
# on numeric-looking strings does not leave a stray NumberFormat behind.
$plainStyle = $ws.Range("D5").Style

$ws.Range("D2").Value = "37.401.83"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "2.063.61"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'237.33"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("E6").Value = "  +3.45%  "
$ws.Range("D7").Value = "'58.35"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("D10").Value = "'57.89"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("D11").Value = "'0.0763"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'0.101"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("D13").Value = "2.366.20"
$ws.Range("E13").Value = "  +4.66%  "
$ws.Range("D14").Value = "'14.46"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  +4.54%  "
$ws.Range("D15").Value = "'21.16"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  +6.55%  "
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("D18").Value = "2.114.31"
$ws.Range("E18").Value = "  +7.33%  "
$ws.Range("D19").Value = "37.512.73"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("E20").Value = "  +18.87%  "
$ws.Range("D21").Value = "'69.11"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "'225.41"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +5.48%  "
$ws.Range("D26").Value = "'2.40"
$ws.Range("D26").Style = $plainStyle
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").Value = "'163.88"
$ws.Range("D27").Style = $plainStyle
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("E28").Value = "  +10.74%  "
$ws.Range("E29").Value = "  +4.71%  "
$ws.Range("D30").Value = "'19.17"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  +2.88%  "
$ws.Range("D31").Value = "'0.126"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "  +2.44%  "
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "'4.51"
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +4.34%  "
$ws.Range("E35").Value = "  +12.45%  "
$ws.Range("D36").Value = "'4.48"
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = "  +6.03%  "
$ws.Range("D37").Value = "'3.38"
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = "  +5.86%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'1.79"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  +13.84%  "
$ws.Range("E41").Value = "  +11.91%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "'4.50"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +25.36%  "
$ws.Range("D44").Value = "'97.63"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  +11.73%  "
$ws.Range("D45").Value = "1.480.63"
$ws.Range("E45").Value = "  +3.54%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'1.16"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +7.58%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0210"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("E48").Value = "  +8.75%  "
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +7.26%  "
$ws.Range("D51").Value = "'2.94"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "  +2.95%  "
